$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("liste_naz")

# Header for new column J
$ws.Range("J1").Value = "AREA"

# Political "AREA" classification per row
$ws.Range("J2").Value  = "DX"
$ws.Range("J3").Value  = "CENTRO"
$ws.Range("J4").Value  = "DX"
$ws.Range("J5").Value  = "DX"
$ws.Range("J6").Value  = "M5S"
$ws.Range("J7").Value  = "CENTRO"
$ws.Range("J8").Value  = "CENTRO"
$ws.Range("J9").Value  = "DX"
$ws.Range("J10").Value = "SX-VERDI"

# Update selected cell to match the final state recorded in the file
$ws.Range("J11").Select()
